# Swap the "Enterprises (absolute #)" row and the
# "Enterprises density (per 1000 people)" row so that the density row
# now comes first (row 12) and the absolute-count row comes second (row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("D12").Value = "'18.5"

$ws.Range("A13").Value = "Enterprises (absolute #)"
$ws.Range("D13").Value = "'34403"
